$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This edit adds a brand-new "4D Box" result row for 25/6/2025 (Wed) at the
# top of the results table (row 2), pushing the previously-listed entries
# (22/6, 21/6, 18/6, 15/6) down by one row. Row 6, which used to be a blank
# placeholder, now receives the 15/6 entry. We also grow the trailing block
# of blank placeholder rows by one (new blank row 35) and backfill a missing
# blank C9 cell so column C's placeholder block stays in sync with column B.

# 1) Shift the four existing data rows (2-5) down into rows 3-6, working
#    from the bottom up so we never clobber a row before reading it.
for ($r = 5; $r -ge 2; $r--) {
    $dest = $r + 1
    $ws.Range("A$dest").Value2 = $ws.Range("A$r").Value2
    $ws.Range("B$dest").Value2 = $ws.Range("B$r").Value2
    $ws.Range("C$dest").Value2 = $ws.Range("C$r").Value2
}

# 2) Write the brand-new top entry into row 2.
$ws.Range("A2").Value2 = "25/6/2025 (Wed)"
$ws.Range("B2").Value2 = "2 0 5 7`n6 2 8 8`n9 1 3 0`n7 6 1 4"
$ws.Range("C2").Value2 = "✅ Direct: 12/4302 (0.28%)`n✅ iBet: 12/226 (5.31%)"

# Row 2 previously had no explicit row height; writing multi-line wrapped
# text into it makes the host auto-stamp a custom height, so re-fit it to
# drop back to the sheet's default (no customHeight override), matching the
# other freshly-written header-adjacent rows.
$ws.Rows.Item(2).AutoFit()

# 3) The placeholder rows below the data keep column B (and sometimes C)
#    pre-formatted but empty. Row 9 only had B9; give it a matching, empty,
#    styled C9 cell (mirrors rows 6-8 which already have both).
$ws.Range("B9").Copy($ws.Range("C9"))

# 4) Extend the placeholder block by one more blank, styled row (35),
#    copied from the last existing placeholder row (34).
$ws.Range("B34").Copy($ws.Range("B35"))

Write-Output "4D box table updated with 25/6/2025 (Wed) results"
